$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reference number format for date column (column A), matching existing rows (style index 2 / YYYY-MM-DD)
$dateFormat = $ws.Range("A2").NumberFormat

# Row 505
$ws.Range("A505").Value = 46024
$ws.Range("A505").NumberFormat = $dateFormat
$ws.Range("B505").Value = "Washington Wizards"
$ws.Range("C505").Value = "Brooklyn Nets"
$ws.Range("D505").Value = -2.5
$ws.Range("E505").Value = 119
$ws.Range("F505").Value = 99
$ws.Range("G505").Value = 17.5

# Row 506
$ws.Range("A506").Value = 46024
$ws.Range("A506").NumberFormat = $dateFormat
$ws.Range("B506").Value = "Indiana Pacers"
$ws.Range("C506").Value = "San Antonio Spurs"
$ws.Range("D506").Value = 4.5
$ws.Range("E506").Value = 113
$ws.Range("F506").Value = 123
$ws.Range("G506").Value = -5.5

# Row 507
$ws.Range("A507").Value = 46024
$ws.Range("A507").NumberFormat = $dateFormat
$ws.Range("B507").Value = "New York Knicks"
$ws.Range("C507").Value = "Atlanta Hawks"
$ws.Range("D507").Value = -4.5
$ws.Range("E507").Value = 99
$ws.Range("F507").Value = 111
$ws.Range("G507").Value = -16.5

# Row 508
$ws.Range("A508").Value = 46024
$ws.Range("A508").NumberFormat = $dateFormat
$ws.Range("B508").Value = "Cleveland Cavaliers"
$ws.Range("C508").Value = "Denver Nuggets"
$ws.Range("D508").Value = -13.5
$ws.Range("E508").Value = 113
$ws.Range("F508").Value = 108
$ws.Range("G508").Value = -8.5

# Row 509
$ws.Range("A509").Value = 46024
$ws.Range("A509").NumberFormat = $dateFormat
$ws.Range("B509").Value = "Chicago Bulls"
$ws.Range("C509").Value = "Orlando Magic"
$ws.Range("D509").Value = 5.5
$ws.Range("E509").Value = 121
$ws.Range("F509").Value = 114
$ws.Range("G509").Value = 12.5

# Row 510
$ws.Range("A510").Value = 46024
$ws.Range("A510").NumberFormat = $dateFormat
$ws.Range("B510").Value = "Milwaukee Bucks"
$ws.Range("C510").Value = "Charlotte Hornets"
$ws.Range("D510").Value = -5.5
$ws.Range("E510").Value = 122
$ws.Range("F510").Value = 121
$ws.Range("G510").Value = -4.5

# Row 511
$ws.Range("A511").Value = 46024
$ws.Range("A511").NumberFormat = $dateFormat
$ws.Range("B511").Value = "New Orleans Pelicans"
$ws.Range("C511").Value = "Portland Trail Blazers"
$ws.Range("D511").Value = 3.5
$ws.Range("E511").Value = 109
$ws.Range("F511").Value = 122
$ws.Range("G511").Value = -9.5

# Row 512
$ws.Range("A512").Value = 46024
$ws.Range("A512").NumberFormat = $dateFormat
$ws.Range("B512").Value = "Phoenix Suns"
$ws.Range("C512").Value = "Sacramento Kings"
$ws.Range("D512").Value = -13.5
$ws.Range("E512").Value = 129
$ws.Range("F512").Value = 102
$ws.Range("G512").Value = 13.5

# Row 513
$ws.Range("A513").Value = 46024
$ws.Range("A513").NumberFormat = $dateFormat
$ws.Range("B513").Value = "Golden State Warriors"
$ws.Range("C513").Value = "Oklahoma City Thunder"
$ws.Range("D513").Value = 13.5
$ws.Range("E513").Value = 94
$ws.Range("F513").Value = 131
$ws.Range("G513").Value = -23.5

# Row 514
$ws.Range("A514").Value = 46024
$ws.Range("A514").NumberFormat = $dateFormat
$ws.Range("B514").Value = "Los Angeles Lakers"
$ws.Range("C514").Value = "Memphis Grizzlies"
$ws.Range("D514").Value = -5.5
$ws.Range("E514").Value = 128
$ws.Range("F514").Value = 121
$ws.Range("G514").Value = 1.5

# Row 515
$ws.Range("A515").Value = 46025
$ws.Range("A515").NumberFormat = $dateFormat
$ws.Range("B515").Value = "Miami Heat"
$ws.Range("C515").Value = "Minnesota Timberwolves"
$ws.Range("D515").Value = 2.5
$ws.Range("E515").Value = 115
$ws.Range("F515").Value = 125
$ws.Range("G515").Value = -7.5

# Row 516
$ws.Range("A516").Value = 46025
$ws.Range("A516").NumberFormat = $dateFormat
$ws.Range("B516").Value = "New York Knicks"
$ws.Range("C516").Value = "Philadelphia 76ers"
$ws.Range("D516").Value = -3.5
$ws.Range("E516").Value = 119
$ws.Range("F516").Value = 130
$ws.Range("G516").Value = -14.5

# Row 517
$ws.Range("A517").Value = 46025
$ws.Range("A517").NumberFormat = $dateFormat
$ws.Range("B517").Value = "Toronto Raptors"
$ws.Range("C517").Value = "Atlanta Hawks"
$ws.Range("D517").Value = -5.5
$ws.Range("E517").Value = 134
$ws.Range("F517").Value = 117
$ws.Range("G517").Value = 11.5

# Row 518
$ws.Range("A518").Value = 46025
$ws.Range("A518").NumberFormat = $dateFormat
$ws.Range("B518").Value = "Chicago Bulls"
$ws.Range("C518").Value = "Charlotte Hornets"
$ws.Range("D518").Value = -2.5
$ws.Range("E518").Value = 99
$ws.Range("F518").Value = 112
$ws.Range("G518").Value = -15.5

# Row 519
$ws.Range("A519").Value = 46025
$ws.Range("A519").NumberFormat = $dateFormat
$ws.Range("B519").Value = "San Antonio Spurs"
$ws.Range("C519").Value = "Portland Trail Blazers"
$ws.Range("D519").Value = -8.5
$ws.Range("E519").Value = 110
$ws.Range("F519").Value = 115
$ws.Range("G519").Value = -13.5

# Row 520
$ws.Range("A520").Value = 46025
$ws.Range("A520").NumberFormat = $dateFormat
$ws.Range("B520").Value = "Dallas Mavericks"
$ws.Range("C520").Value = "Houston Rockets"
$ws.Range("D520").Value = 7.5
$ws.Range("E520").Value = 110
$ws.Range("F520").Value = 104
$ws.Range("G520").Value = 13.5

# Row 521
$ws.Range("A521").Value = 46025
$ws.Range("A521").NumberFormat = $dateFormat
$ws.Range("B521").Value = "Golden State Warriors"
$ws.Range("C521").Value = "Utah Jazz"
$ws.Range("D521").Value = -12.5
$ws.Range("E521").Value = 123
$ws.Range("F521").Value = 114
$ws.Range("G521").Value = -3.5

# Row 522
$ws.Range("A522").Value = 46025
$ws.Range("A522").NumberFormat = $dateFormat
$ws.Range("B522").Value = "Los Angeles Clippers"
$ws.Range("C522").Value = "Boston Celtics"
$ws.Range("D522").Value = -1.5
$ws.Range("E522").Value = 115
$ws.Range("F522").Value = 146
$ws.Range("G522").Value = -32.5

# Row 523
$ws.Range("A523").Value = 46026
$ws.Range("A523").NumberFormat = $dateFormat
$ws.Range("B523").Value = "Cleveland Cavaliers"
$ws.Range("C523").Value = "Detroit Pistons"
$ws.Range("D523").Value = -3.5
$ws.Range("E523").Value = 110
$ws.Range("F523").Value = 114
$ws.Range("G523").Value = -7.5

# Row 524
$ws.Range("A524").Value = 46026
$ws.Range("A524").NumberFormat = $dateFormat
$ws.Range("B524").Value = "Orlando Magic"
$ws.Range("C524").Value = "Indiana Pacers"
$ws.Range("D524").Value = -6.5
$ws.Range("E524").Value = 135
$ws.Range("F524").Value = 127
$ws.Range("G524").Value = 1.5

# Row 525
$ws.Range("A525").Value = 46026
$ws.Range("A525").NumberFormat = $dateFormat
$ws.Range("B525").Value = "Brooklyn Nets"
$ws.Range("C525").Value = "Denver Nuggets"
$ws.Range("D525").Value = 2.5
$ws.Range("E525").Value = 127
$ws.Range("F525").Value = 115
$ws.Range("G525").Value = 14.5

# Row 526
$ws.Range("A526").Value = 46026
$ws.Range("A526").NumberFormat = $dateFormat
$ws.Range("B526").Value = "Washington Wizards"
$ws.Range("C526").Value = "Minnesota Timberwolves"
$ws.Range("D526").Value = 10.5
$ws.Range("E526").Value = 115
$ws.Range("F526").Value = 141
$ws.Range("G526").Value = -15.5

# Row 527
$ws.Range("A527").Value = 46026
$ws.Range("A527").NumberFormat = $dateFormat
$ws.Range("B527").Value = "Miami Heat"
$ws.Range("C527").Value = "New Orleans Pelicans"
$ws.Range("D527").Value = -7.5
$ws.Range("E527").Value = 125
$ws.Range("F527").Value = 106
$ws.Range("G527").Value = 11.5

# Row 528
$ws.Range("A528").Value = 46026
$ws.Range("A528").NumberFormat = $dateFormat
$ws.Range("B528").Value = "Phoenix Suns"
$ws.Range("C528").Value = "Oklahoma City Thunder"
$ws.Range("D528").Value = 10.5
$ws.Range("E528").Value = 108
$ws.Range("F528").Value = 105
$ws.Range("G528").Value = 13.5

# Row 529
$ws.Range("A529").Value = 46026
$ws.Range("A529").NumberFormat = $dateFormat
$ws.Range("B529").Value = "Sacramento Kings"
$ws.Range("C529").Value = "Milwaukee Bucks"
$ws.Range("D529").Value = 6.5
$ws.Range("E529").Value = 98
$ws.Range("F529").Value = 115
$ws.Range("G529").Value = -10.5

# Row 530
$ws.Range("A530").Value = 46026
$ws.Range("A530").NumberFormat = $dateFormat
$ws.Range("B530").Value = "Los Angeles Lakers"
$ws.Range("C530").Value = "Memphis Grizzlies"
$ws.Range("D530").Value = -5.5
$ws.Range("E530").Value = 120
$ws.Range("F530").Value = 114
$ws.Range("G530").Value = 0.5

# Row 531
$ws.Range("A531").Value = 46027
$ws.Range("A531").NumberFormat = $dateFormat
$ws.Range("B531").Value = "Detroit Pistons"
$ws.Range("C531").Value = "New York Knicks"
$ws.Range("D531").Value = 1.5
$ws.Range("E531").Value = 121
$ws.Range("F531").Value = 90
$ws.Range("G531").Value = 32.5

# Row 532
$ws.Range("A532").Value = 46027
$ws.Range("A532").NumberFormat = $dateFormat
$ws.Range("B532").Value = "Toronto Raptors"
$ws.Range("C532").Value = "Atlanta Hawks"
$ws.Range("D532").Value = -2.5
$ws.Range("E532").Value = 118
$ws.Range("F532").Value = 100
$ws.Range("G532").Value = 15.5

# Row 533
$ws.Range("A533").Value = 46027
$ws.Range("A533").NumberFormat = $dateFormat
$ws.Range("B533").Value = "Boston Celtics"
$ws.Range("C533").Value = "Chicago Bulls"
$ws.Range("D533").Value = -10.5
$ws.Range("E533").Value = 115
$ws.Range("F533").Value = 101
$ws.Range("G533").Value = 3.5

# Row 534
$ws.Range("A534").Value = 46027
$ws.Range("A534").NumberFormat = $dateFormat
$ws.Range("B534").Value = "Oklahoma City Thunder"
$ws.Range("C534").Value = "Charlotte Hornets"
$ws.Range("D534").Value = -15.5
$ws.Range("E534").Value = 97
$ws.Range("F534").Value = 124
$ws.Range("G534").Value = -42.5

# Row 535
$ws.Range("A535").Value = 46027
$ws.Range("A535").NumberFormat = $dateFormat
$ws.Range("B535").Value = "Houston Rockets"
$ws.Range("C535").Value = "Phoenix Suns"
$ws.Range("D535").Value = -8.5
$ws.Range("E535").Value = 100
$ws.Range("F535").Value = 97
$ws.Range("G535").Value = -5.5

# Row 536
$ws.Range("A536").Value = 46027
$ws.Range("A536").NumberFormat = $dateFormat
$ws.Range("B536").Value = "Philadelphia 76ers"
$ws.Range("C536").Value = "Denver Nuggets"
$ws.Range("D536").Value = -14.5
$ws.Range("E536").Value = 124
$ws.Range("F536").Value = 125
$ws.Range("G536").Value = -15.5

# Row 537
$ws.Range("A537").Value = 46027
$ws.Range("A537").NumberFormat = $dateFormat
$ws.Range("B537").Value = "Los Angeles Clippers"
$ws.Range("C537").Value = "Golden State Warriors"
$ws.Range("D537").Value = 3.5
$ws.Range("E537").Value = 103
$ws.Range("F537").Value = 102
$ws.Range("G537").Value = 4.5

# Row 538
$ws.Range("A538").Value = 46027
$ws.Range("A538").NumberFormat = $dateFormat
$ws.Range("B538").Value = "Portland Trail Blazers"
$ws.Range("C538").Value = "Utah Jazz"
$ws.Range("D538").Value = -5.5
$ws.Range("E538").Value = 137
$ws.Range("F538").Value = 117
$ws.Range("G538").Value = 14.5

Write-Output "Added rows 505-538"